# Applies the two substantive changes from the commit:
#  1. Slide 16's table switches to a different (built-in) table style.
#  2. The deck's theme color palette is swapped from the "Integral" custom
#     palette to the standard Office palette (the before/after OOXML shows
#     the theme1.xml / theme2.xml parts exchanging their full contents;
#     the only semantic difference between the two parts, besides their
#     <a:theme>/<a:clrScheme> "name" attributes, is these 10 color values -
#     font & format schemes are byte-for-byte identical already).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$s = $p.Slides.Item(16)
foreach ($sh in $s.Shapes) {
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{BDA18A57-F0C3-4C06-B97D-36D93F9380A4}")
    }
}

# --- 2. Theme colors --------------------------------------------------------
# RGBColor.RGB packs bytes as 0xBBGGRR (classic VBA RGB() long), so the
# literals below are the byte-swapped form of the target RRGGBB hex codes:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6 accent1=5B9BD5
#   accent2=ED7D31 accent3=A5A5A5 accent4=FFC000 accent5=4472C4
#   accent6=70AD47 hlink=0563C1 folHlink=954F72
$targetColors = @(
    0x000000,
    0xFFFFFF,
    0x6A5444,
    0xE6E6E7,
    0xD59B5B,
    0x317DED,
    0xA5A5A5,
    0x00C0FF,
    0xC47244,
    0x47AD70,
    0xC16305,
    0x724F95
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $targetColors[$i - 1]
}
